$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.454.08"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "1.798.16"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.68"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3811"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3460"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.204"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07504"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +8.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.482"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "1.797.74"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.066"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06642"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.87"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.506"
$ws.Range("E21").Value = "  +4.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.36"
$ws.Range("E22").Value = "  +4.23%  "
$ws.Range("D23").Value = "27.426.78"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.52"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.427"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.566"
$ws.Range("E26").Value = "  +5.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.499"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.42"
$ws.Range("E28").Value = "  +9.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.23"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "2.000.58"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.04"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.062"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.136"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08696"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.27"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.689"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.460"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6900"
$ws.Range("E38").Value = "  +10.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.922"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2210"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06364"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02340"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.277"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.39"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6453"
$ws.Range("E45").Value = "  +6.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.865"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.129"
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.00"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07194"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.44"
$ws.Range("E51").Value = "  +2.32%  "
